# fix: correct typo in heading that broke XLS file chunking
#
# The header cell G3 on the "metadata" sheet contained "Description"
# (capitalized) while the code expects the lowercase "description".
# Correct the typo so the heading matches the expected column name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

$ws.Range("G3").Value = "description"

# Update the active selection to match the saved state (cursor moved to G4)
$ws.Range("G4").Select()
